$wb = $excel.ActiveWorkbook

# Sheets that contain the "展览" (exhibition) event rows which need updating:
#   "展览"   - the dedicated exhibition sheet
#   "全部类型" - the combined "all types" sheet mirrors the same rows
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 594   # was 590
    $ws.Range("F6").Value = 42    # was 41
    $ws.Range("F8").Value = 561   # was 552
    $ws.Range("F9").Value = 3723  # was 3711
}
